$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B5 ("Jurek" -> "Kurek") - this corresponds to row 5 second column
$ws.Range("B5").Value = "Kurek"

# Update the active cell selection to D17
$ws.Range("D17").Select()
